# Auto-generated: update TPM-derived values in sheet1 per commit "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6288379999999999
$ws.Range("H2").Value = 1.886514
$ws.Range("I2").Value = 0.003263417952351538
$ws.Range("J2").Value = 0.003263417952351539
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 1.274727151982666
$ws.Range("R2").Value = 11.472544367844
$ws.Range("S2").Value = 0.00002152643346960648
$ws.Range("T2").Value = 0.00002152643346960648
$ws.Range("G3").Value = 0.6288379999999999
$ws.Range("H3").Value = 1.886514
$ws.Range("I3").Value = 0.003263417952351538
$ws.Range("J3").Value = 0.003263417952351539
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 161.2619240772026
$ws.Range("R3").Value = 1451.357316694824
$ws.Range("S3").Value = 0.002723244793546091
$ws.Range("T3").Value = 0.002723244793546091
$ws.Range("G4").Value = 0.6288379999999999
$ws.Range("H4").Value = 1.886514
$ws.Range("I4").Value = 0.003263417952351538
$ws.Range("J4").Value = 0.003263417952351539
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 30.71261498129533
$ws.Range("R4").Value = 276.413534831658
$ws.Range("S4").Value = 0.0005186467253358412
$ws.Range("T4").Value = 0.0005186467253358411
$ws.Range("I5").Value = 0.8672450858064795
$ws.Range("J5").Value = 0.8672450858064795
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 338.7555239452122
$ws.Range("R5").Value = 3048.79971550691
$ws.Range("S5").Value = 0.005720595373940424
$ws.Range("T5").Value = 0.005720595373940423
$ws.Range("I6").Value = 0.8672450858064795
$ws.Range("J6").Value = 0.8672450858064795
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.7236954319470881
$ws.Range("T6").Value = 0.723695431947088
$ws.Range("I7").Value = 0.8672450858064795
$ws.Range("J7").Value = 0.8672450858064795
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 8161.799929917667
$ws.Range("R7").Value = 73456.199369259
$ws.Range("S7").Value = 0.1378290584854511
$ws.Range("T7").Value = 0.137829058485451
$ws.Range("G8").Value = 24.952113
$ws.Range("H8").Value = 74.85633900000001
$ws.Range("I8").Value = 0.129491496241169
$ws.Range("J8").Value = 0.129491496241169
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 50.580810861366
$ws.Range("R8").Value = 455.227297752294
$ws.Range("S8").Value = 0.0008541627580085855
$ws.Range("T8").Value = 0.0008541627580085856
$ws.Range("G9").Value = 24.952113
$ws.Range("H9").Value = 74.85633900000001
$ws.Range("I9").Value = 0.129491496241169
$ws.Range("J9").Value = 0.129491496241169
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 6398.827284883836
$ws.Range("R9").Value = 57589.44556395453
$ws.Range("S9").Value = 0.1080575789237033
$ws.Range("T9").Value = 0.1080575789237033
$ws.Range("G10").Value = 24.952113
$ws.Range("H10").Value = 74.85633900000001
$ws.Range("I10").Value = 0.129491496241169
$ws.Range("J10").Value = 0.129491496241169
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 1218.667827864687
$ws.Range("R10").Value = 10968.01045078219
$ws.Range("S10").Value = 0.02057975455945708
$ws.Range("T10").Value = 0.02057975455945708

Write-Host "Updated" 94 "cells with new TPM values."
